$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price observation was inserted as row 181 (La Araucania /
# Vega Modelo de Temuco / Ciboulette), pushing the existing rows 181:197
# down to 182:198 and extending the used range to A1:R198.
$ws.Rows.Item(181).Insert()

# Populate the newly inserted row 181 with the new observation. All the
# "constant" columns for this block (A,B,C,E,F,G,H,I,N,O,Q,R) keep the same
# values used throughout this subset.
$ws.Cells.Item(181, 1).Value = 10
$ws.Cells.Item(181, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(181, 3).Value = "La Araucanía"

$ws.Cells.Item(181, 4).Value = (Get-Date -Year 2022 -Month 1 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Cells.Item(181, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item(181, 5).Value = 9
$ws.Cells.Item(181, 6).Value = 100112039
$ws.Cells.Item(181, 7).Value = "Ciboulette"
$ws.Cells.Item(181, 8).Value = "Sin especificar"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 50
$ws.Cells.Item(181, 11).Value = 5000
$ws.Cells.Item(181, 12).Value = 5000
$ws.Cells.Item(181, 13).Value = 5000
$ws.Cells.Item(181, 14).Value = "`$/docena de atados"
$ws.Cells.Item(181, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(181, 16).Value = 1667
$ws.Cells.Item(181, 17).Value = 3
$ws.Cells.Item(181, 18).Value = "Hortaliza"
